# Update Capex logic to be an initial capital outlay instead of annual deduction
#
# Concretely (per the target OOXML diff) this commit reshapes the tail of the
# deck:
#   - Slide 6 "Financial Overview & Sensitivities" loses its sensitivity
#     analysis paragraphs and is retitled "Financial Overview".
#   - A brand-new "Sensitivities" slide (title + {{SENSITIVITY_ANALYSIS}}
#     placeholder) is inserted right after it, becoming slide 7.
#   - The former slide 7 ("Appendix") shifts down to become the new slide 8,
#     unchanged in content.

# Helper: replace a TextRange's text without PowerPoint stamping a fresh
# <a:rPr lang="en-US"/> on the (only) run, which happens when assigning
# .Text directly on a single-paragraph range. Inserting the new text before
# the old one (extending the existing run) and then trimming the old
# characters off the tail keeps the XML minimal/clean, matching how the
# original template's runs look (no rPr).
function Set-CleanText {
    param($TextRange, [string]$NewText)

    $oldLen = $TextRange.Length
    $TextRange.InsertBefore($NewText)
    if ($oldLen -gt 0) {
        $newLen = $TextRange.Length
        $staleTail = $TextRange.Characters($newLen - $oldLen + 1, $oldLen)
        $staleTail.Delete()
    }
}

$p = $ppt.ActivePresentation

# --- 1. Slide 6: "Financial Overview & Sensitivities" -> "Financial Overview"
#     and drop the trailing "Sensitivity Analysis:" / {{SENSITIVITY_ANALYSIS}}
#     paragraphs (they move to their own slide below).
$finSlide = $p.Slides.Item(6)
Set-CleanText $finSlide.Shapes.Item(1).TextFrame.TextRange "Financial Overview"
$finSlide.Shapes.Item(2).TextFrame.TextRange.Text = `
    "Key Metrics:`rEntry Yield: {{ENTRY_YIELD}}`rIRR: {{IRR}}`rEquity Multiple: {{MOIC}}`rExit Yield: {{EXIT_YIELD}}"

# --- 2. Duplicate the current slide 7 ("Appendix"). The duplicate is placed
#     right after it (new slide 8) and is byte-for-byte identical (title,
#     body, layout, color-map override, etc.) -- exactly the "Appendix slide
#     moves to the end" half of the change.
$appendixSlide = $p.Slides.Item(7)
$appendixSlide.Duplicate() | Out-Null

# --- 3. Turn the original slide 7 into the new "Sensitivities" slide, reusing
#     its Title/Content placeholders.
Set-CleanText $appendixSlide.Shapes.Item(1).TextFrame.TextRange "Sensitivities"
Set-CleanText $appendixSlide.Shapes.Item(2).TextFrame.TextRange "{{SENSITIVITY_ANALYSIS}}"
